$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Microstate IDs whose (ID, SMILES) pair must be removed from the table.
# (resonance structures / duplicate geometric isomers eliminated per v1.3.2)
$removed = @{}
$removed["SM18_micro006"] = $true
$removed["SM18_micro012"] = $true
$removed["SM18_micro014"] = $true
$removed["SM18_micro018"] = $true
$removed["SM18_micro022"] = $true
$removed["SM18_micro023"] = $true
$removed["SM18_micro024"] = $true
$removed["SM18_micro036"] = $true
$removed["SM18_micro051"] = $true

# Collect the existing (ID, SMILES) pairs from the data rows (3..76),
# keeping only the ones that were not eliminated.
$ids = @()
$smiles = @()
for ($r = 3; $r -le 76; $r++) {
  $idAddr = "B" + $r
  $smiAddr = "C" + $r
  $id = $ws.Range($idAddr).Text
  $smi = $ws.Range($smiAddr).Text
  if (-not $removed.ContainsKey($id)) {
    $ids += $id
    $smiles += $smi
  }
}

# Write the surviving pairs back into rows 3..67 (65 rows), preserving
# each row's existing style/formatting - only .Value is touched.
for ($i = 0; $i -lt $ids.Count; $i++) {
  $r = $i + 3
  $ws.Range("B" + $r).Value = $ids[$i]
  $ws.Range("C" + $r).Value = $smiles[$i]
}

# Drop the now-unused trailing rows (68..76).
$ws.Range("A68:D76").EntireRow.Delete()

# Drop the trailing 9 microstate depiction pictures (the tail of the
# Shapes collection) to match the shrunk table.
$shapeCount = $ws.Shapes.Count
$firstToDelete = $shapeCount - 8
for ($i = $shapeCount; $i -ge $firstToDelete; $i--) {
  $ws.Shapes.Item($i).Delete()
}

# Update the sheet title / header cell.
$ws.Range("A1").Value = "Microstate List"
